$wb = $excel.ActiveWorkbook

# --- Sheet 1 (physical sheet1.xml): rename -> "NB_TO-1651589045585127" ---
# Grows from 5 rows (A1:B5) to 10 rows (A1:B10): add new rows 6-10, copying
# the existing style (col A, style s="1") from row 2 down into the new rows.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "NB_TO-1651589045585127"

$ws1.Range("A2").Copy()
$ws1.Range("A6:A10").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "ZB-match_5-1651589043832483.csv"
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "OB-16515890441159182.csv"
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "OB-16515890448441927.csv"
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "TB-16515890454109972.csv"
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "TB-1651589044953564.csv"
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "ZB-match_9-16515890438481045.csv"
$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "TB-16515890455694811.csv"
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "OB-16515890447639709.csv"
$ws1.Range("A10").Value = 8
$ws1.Range("B10").Value = "ZB-match_3-16515890438012626.csv"

# --- Sheet 2 (physical sheet2.xml): rename -> "TOL_TO-16515890456319466" ---
# Shrinks from 10 rows (A1:B10) to 7 rows (A1:B7): delete old rows 8-10.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TOL_TO-16515890456319466"

$ws2.Range("A8:B10").EntireRow.Delete()

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "MM_stims-16515890456006975.csv"
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = "ZM_stims-1651589045585127.csv"
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "MM_stims-16515890456163235.csv"
$ws2.Range("A5").Value = 3
$ws2.Range("B5").Value = "ZM_stims-16515890456006975.csv"
$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "MM_stims-16515890456319466.csv"
$ws2.Range("A7").Value = 5
$ws2.Range("B7").Value = "ZM_stims-16515890456163235.csv"

# --- Sheet 3 (physical sheet3.xml): rename -> "RS_TO-16515890456319466" ---
# Stays 3 rows (A1:B3); swap "eyes open" / "eyes closed".
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16515890456319466"

$ws3.Range("A2").Value = 0
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (physical sheet4.xml): rename -> "GNG_TO-16515890456631956" ---
# Shrinks from 7 rows (A1:B7) to 5 rows (A1:B5): delete old rows 6-7.
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "GNG_TO-16515890456631956"

$ws4.Range("A6:B7").EntireRow.Delete()

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "go_stims-16515890456319466.csv"
$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "GNG_stims-16515890456475704.csv"
$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = "go_stims-16515890456475704.csv"
$ws4.Range("A5").Value = 3
$ws4.Range("B5").Value = "GNG_stims-16515890456631956.csv"

# --- Sheet 5 (physical sheet5.xml): rename -> "vSAT_TO-16515890457256973" ---
# Stays 5 rows (A1:B5).
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16515890457256973"

$ws5.Range("A2").Value = 0
$ws5.Range("B2").Value = "vSAT_stims-16515890457100737.csv"
$ws5.Range("A3").Value = 1
$ws5.Range("B3").Value = "SAT_stims-16515890456631956.csv"
$ws5.Range("A4").Value = 2
$ws5.Range("B4").Value = "vSAT_stims-16515890456944494.csv"
$ws5.Range("A5").Value = 3
$ws5.Range("B5").Value = "SAT_stims-16515890456788225.csv"
